$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("BF2")
$c.Formula = "'2008-02-15"
Write-Output "VAL2:$($c.Value2)"
Write-Output "TEXT:$($c.Text)"
